$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- USERS table (rows 1-7) ---
# Row 2 used to hold the "UserPk" surrogate-key field row; it is now cleared
# and given a blank/white fill (matches the new HOME table treatment below).
$ws.Range("A2:E2").ClearContents()
$ws.Range("A2:E2").Interior.ThemeColor = 2

# Row 4 ("email" field) is promoted to the table's primary key: highlight it
# with the same green fill used for the other PK rows.
$ws.Range("A4:E4").Interior.Color = $ws.Range("A21").Interior.Color

# --- HOME table (rows 9-18) ---
# Row 10 used to hold the "HomePk" surrogate-key field row; clear + blank fill.
$ws.Range("A10:E10").ClearContents()
$ws.Range("A10:E10").Interior.ThemeColor = 2

# Row 12 ("Home_Name" field) becomes the new primary key: was UNIQUE, now
# PRIMARY KEY, highlighted like the other PK rows.
$ws.Range("D12").Value = "PRIMARY KEY"
$ws.Range("A12:E12").Interior.Color = $ws.Range("A21").Interior.Color

# --- USERACTIVITY bridge table (rows 44-48) ---
# The combined-primary-key row used to reference "UserPk"/INT; now it
# references the new "email"/varchar (60) primary key column.
$ws.Range("A45").Value = "email"
$ws.Range("B45").Value = "varchar (60)"

# --- View state ---
$ws.Range("C9").Select()
